# Set the "Remontable" column (K) to 1 for the rows where it is currently 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,41,45,46,47,50,51,52,53,54,55,56)

foreach ($r in $rows) {
    $ws.Range("K$r").Value = 1
}
